$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("E21").Formula = "=1000*0.000196236724150367"
$ws.Range("F21").Formula = "=1000*0.00177883368451148"
$ws.Range("G21").Formula = "=1000*0.000814738101325929"
$ws.Range("H21").Formula = "=1000*0.00127509713638574"
$ws.Range("I21").Formula = "=1000*0.000012051522389811"
$ws.Range("J21").Formula = "=1000*0.00120656989626586"

$ws.Range("B17").Select()
